$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "edit1"
$ws.Range("B11").Value = "riya-morankar"
$ws.Range("C11").Value = "Squashed"

# "2025-06-18" must land as literal text (matching the other Date-column
# cells), not get auto-parsed into a date serial number. Force the cell to
# Text format before assigning, then drop back to the workbook's default
# style so no stray number-format attribute is left on the cell.
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2025-06-18"
$ws.Range("E11").Style = "Normal"

$ws.Range("F11").Value = "N/A"
